$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data, sorted descending by value, with Swedish and Uzbek removed.
$data = @(
    @("English", 21.44696729049518),
    @("Chinese", 19.49450511465953),
    @("Spanish", 6.16332172048769),
    @("German", 4.268509755873809),
    @("Arabic", 4.16374108609093),
    @("Japanese", 3.987400586714449),
    @("Russian", 3.235507227495868),
    @("Malay-Indonesian", 3.157883726556382),
    @("Portuguese", 2.836017439400714),
    @("French", 2.505005038123243),
    @("Turkish", 1.911955492510335),
    @("Italian", 1.84879646833128),
    @("Korean", 1.741991329786734),
    @("Dutch", 1.209315284991469),
    @("Persian", 1.058255940709667),
    @("Polish", 0.9797889098155134),
    @("Thai", 0.9535397572350784),
    @("Urdu", 0.9082800637918591),
    @("Vietnamese", 0.7997661210102581),
    @("Bengali", 0.7899811735361185)
)

# Clear out the old data area (rows 2-23, columns A:B) first.
$ws.Range("A2:B23").Clear()

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Re-apply the same cell format (bold, centered/top aligned, thin box border)
# used by the header/label column to the new column A label cells.
$ws.Range("A1").Copy()
$ws.Range("A2:A21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wb.Save()
